$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: preserve old row 16 (Gaussian-Quadrature) into a scratch row, since its data moves to row 10 ---
# (bounded A:M range used instead of whole-row copy so the sheet dimension does not bleed out to column XFD)
$ws.Range("A16:M16").Copy($ws.Range("A100:M100"))

# --- Step 2: shift rows 10-15 down to 14-19 (copy in descending src order so we never read an already-overwritten row) ---
$ws.Range("A15:M15").Copy($ws.Range("A19:M19"))
$ws.Range("A14:M14").Copy($ws.Range("A18:M18"))
$ws.Range("A13:M13").Copy($ws.Range("A17:M17"))
$ws.Range("A12:M12").Copy($ws.Range("A16:M16"))
$ws.Range("A11:M11").Copy($ws.Range("A15:M15"))
$ws.Range("A10:M10").Copy($ws.Range("A14:M14"))

# --- Step 3: move the preserved Gaussian-Quadrature row into its new position (row 10), then wipe the scratch row ---
$ws.Range("A100:M100").Copy($ws.Range("A10:M10"))
$ws.Rows(100).Clear()

# --- Step 4: add the 3 new Spiral rows (11-13); formatting already correct (untouched original rows) ---
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 1.17724958549151
$ws.Range("D11").Value = 0.4610126382763439
$ws.Range("E11").Value = 1.022395423802949
$ws.Range("F11").Value = 1.17724958549151
$ws.Range("G11").Value = 0.8914036606228333
$ws.Range("H11").Value = 0.8926078569125404
$ws.Range("I11").Value = 1.068184896857687
$ws.Range("J11").Value = 0.4610126382763439
$ws.Range("K11").Value = 0.7417040310396463
$ws.Range("L11").Value = 0.9594768082655782
$ws.Range("M11").Value = 0.9188090103273107

$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 1.17197013134176
$ws.Range("D12").Value = 0.4619300435693042
$ws.Range("E12").Value = 1.024057898114124
$ws.Range("F12").Value = 1.17197013134176
$ws.Range("G12").Value = 0.8924799158639082
$ws.Range("H12").Value = 0.8930742401140037
$ws.Range("I12").Value = 1.068063741290045
$ws.Range("J12").Value = 0.4619300435693042
$ws.Range("K12").Value = 0.7429939708417144
$ws.Range("L12").Value = 0.9574820510917372
$ws.Range("M12").Value = 0.9185959950488577

$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 1.17687957816677
$ws.Range("D13").Value = 0.4609012309710387
$ws.Range("E13").Value = 1.022863514736807
$ws.Range("F13").Value = 1.17687957816677
$ws.Range("G13").Value = 0.8913523941060147
$ws.Range("H13").Value = 0.891445048164721
$ws.Range("I13").Value = 1.068198378089818
$ws.Range("J13").Value = 0.4609012309710387
$ws.Range("K13").Value = 0.7418823728539231
$ws.Range("L13").Value = 0.9593809755103467
$ws.Range("M13").Value = 0.9186066907058618

# --- Step 5: refresh column A running index for every data row (1..17) ---
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12
$ws.Range("A15").Value = 13
$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15
$ws.Range("A18").Value = 16
$ws.Range("A19").Value = 17

# --- Step 6: row 17 (HexGrid-90degTilt5degRes) had its C/F recomputed with a tiny float nudge ---
$ws.Range("C17").Value = 0.9963614538049353
$ws.Range("F17").Value = 0.9963614538049353
